$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
    $ws.Range($cellRef).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = $false
}

Set-TextValue 'D2' '38.159.50'
Set-TextValue 'E2' '  +3.38%  '
Set-TextValue 'D3' '2.058.89'
Set-TextValue 'E3' '  +3.33%  '
Set-TextValue 'E4' '  -0.50%  '
Set-TextValue 'D5' '230.27'
Set-TextValue 'E5' '  +2.86%  '
Set-TextValue 'E6' '  +2.01%  '
Set-TextValue 'D7' '58.42'
Set-TextValue 'E7' '  +7.95%  '
Set-TextValue 'E8' '  -0.04%  '
Set-TextValue 'D9' '0.388'
Set-TextValue 'E9' '  +3.64%  '
Set-TextValue 'D10' '0.0809'
Set-TextValue 'E10' '  +4.13%  '
Set-TextValue 'E11' '  +0.39%  '
Set-TextValue 'D12' '2.362.88'
Set-TextValue 'E12' '  +3.23%  '
Set-TextValue 'D13' '14.63'
Set-TextValue 'E13' '  +4.61%  '
Set-TextValue 'D14' '20.65'
Set-TextValue 'E14' '  +3.67%  '
Set-TextValue 'E15' '  +3.01%  '
Set-TextValue 'D16' '5.29'
Set-TextValue 'E16' '  +4.51%  '
Set-TextValue 'D17' '2.056.08'
Set-TextValue 'E17' '  +3.16%  '
Set-TextValue 'D18' '38.053.72'
Set-TextValue 'E18' '  +3.35%  '
Set-TextValue 'E19' '  +1.74%  '
Set-TextValue 'D20' '69.91'
Set-TextValue 'E20' '  +2.36%  '
Set-TextValue 'E21' '  +3.15%  '
Set-TextValue 'D22' '224.83'
Set-TextValue 'E22' '  +1.49%  '
Set-TextValue 'E23' '  -0.01%  '
Set-TextValue 'E24' '  +0.97%  '
Set-TextValue 'E25' '  +4.72%  '
Set-TextValue 'D26' '9.32'
Set-TextValue 'E26' '  +2.90%  '
Set-TextValue 'D27' '166.28'
Set-TextValue 'E28' '  +8.98%  '
Set-TextValue 'D29' '19.06'
Set-TextValue 'E29' '  +2.81%  '
Set-TextValue 'E30' '  +2.08%  '
Set-TextValue 'E31' '  +2.64%  '
Set-TextValue 'E32' '  +2.09%  '
Set-TextValue 'E33' '  +5.96%  '
Set-TextValue 'E34' '  +1.55%  '
Set-TextValue 'E35' '  +7.58%  '
Set-TextValue 'E36' '  +1.70%  '
Set-TextValue 'D37' '6.04'
Set-TextValue 'E37' '  +16.38%  '
Set-TextValue 'E38' '  +6.63%  '
Set-TextValue 'E39' '  -0.08%  '
Set-TextValue 'D40' '98.42'
Set-TextValue 'E40' '  +4.66%  '
Set-TextValue 'E41' '  +2.82%  '
Set-TextValue 'B42' 'Maker'
Set-TextValue 'C42' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D42' '1.481.15'
Set-TextValue 'E42' '  +1.68%  '
Set-TextValue 'B43' 'InjectiveProtocol'
Set-TextValue 'C43' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D43' '16.91'
Set-TextValue 'E43' '  +4.53%  '
Set-TextValue 'D44' '0.0946'
Set-TextValue 'E44' '  +3.65%  '
Set-TextValue 'D45' '2.85'
Set-TextValue 'E45' '  +4.39%  '
Set-TextValue 'E46' '  +1.00%  '
Set-TextValue 'D47' '4.10'
Set-TextValue 'E47' '  +18.45%  '
Set-TextValue 'E48' '  +2.20%  '
Set-TextValue 'E49' '  +2.61%  '
Set-TextValue 'D50' '7.08'
Set-TextValue 'E50' '  -0.36%  '
Set-TextValue 'D51' '2.250.92'
Set-TextValue 'E51' '  +3.17%  '
